$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '330.01'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '1.36%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '44.16'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-0.94%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.502'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-1.81%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.08003'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-0.70%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.976'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '3.87%'
$ws.Range("B7").Value = 'GateToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '4.391'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '2.02%'
$ws.Range("B8").Value = 'BTSEToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '2.574'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-4.68%'
$ws.Range("B9").Value = 'MXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.9516'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '0.97%'
$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1119'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-3.93%'
$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1890'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '1.77%'
$ws.Range("B12").Value = 'MCDex'
$ws.Range("C12").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '10.56'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '25.25%'
$ws.Range("B13").Value = 'MandalaExchangeToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09914'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '0.13%'
$ws.Range("B14").Value = 'BitrueCoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.04785'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '12.07%'
$ws.Range("B15").Value = 'BitMartToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.1065'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-0.13%'
$ws.Range("B16").Value = 'BitForexToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.001267'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-1.17%'
$ws.Range("B17").Value = 'CoinExToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.04080'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-3.10%'
$ws.Range("B18").Value = 'TigerCash'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.005996'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '0.23%'
$ws.Range("B19").Value = 'HotbitToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.004375'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-2.06%'
$ws.Range("B20").Value = 'LEO'
$ws.Range("C20").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.369'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '-6.23%'
$ws.Range("B21").Value = 'BitpandaEcosystemToken'
$ws.Range("C21").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.3493'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-0.16%'
$ws.Range("B22").Value = 'ProBitToken'
$ws.Range("C22").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.1417'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '3.31%'
$ws.Range("B23").Value = 'ZBToken'
$ws.Range("C23").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.2587'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-0.98%'
$ws.Range("B24").Value = 'BitKan'
$ws.Range("C24").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001270'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '2.10%'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '1.57%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0003746'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-6.26%'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02591'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '-0.87%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05696'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '4.90%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.007560'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-2.01%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1400'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '0.36%'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '3.05%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002016'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-0.54%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.008359'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-2.50%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00007135'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-0.03%'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-0.12%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0005802'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-0.15%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.003532'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '55.31%'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-3.17%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002101'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.12%'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0002001'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '-0.12%'
